$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 19517
$ws1.Range("F8").Value = 7
$ws1.Range("F15").Value = 0
$ws1.Range("F18").Value = 1327
$ws1.Range("F19").Value = 0
$ws1.Range("F24").Value = 0
$ws1.Range("F26").Value = 1067
$ws1.Range("F29").Value = 164
$ws1.Range("F30").Value = 5220
$ws1.Range("F31").Value = 552
$ws1.Range("F32").Value = 47
$ws1.Range("F33").Value = 1624
$ws1.Range("F37").Value = 12487
$ws1.Range("F44").Value = 3975

# Sheet "全部类型" (all types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 19517
$ws4.Range("F8").Value = 7
$ws4.Range("F12").Value = 249
$ws4.Range("F19").Value = 363
$ws4.Range("F25").Value = 0
$ws4.Range("F26").Value = 1067
$ws4.Range("F29").Value = 164
$ws4.Range("F30").Value = 5220
$ws4.Range("F31").Value = 552
$ws4.Range("F33").Value = 47
$ws4.Range("F35").Value = 1627
$ws4.Range("F36").Value = 24
$ws4.Range("F39").Value = 12487
$ws4.Range("F42").Value = 10
$ws4.Range("F46").Value = 3975
